$p = $ppt.ActivePresentation

# The deck currently ends with slide 24 ("Thank You"). Insert a new
# "GitHub Link" slide right before it (i.e. at position 24), using the
# same "Title and Content" layout (slideLayout2.xml / CustomLayout 2)
# that slide already uses. This pushes "Thank You" to position 25.
$newSlide = $p.Slides.Add(24, [PowerPoint.PpSlideLayout]::ppLayoutText)

# --- Title placeholder ---
$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "GitHub Link"

# --- Content placeholder ---
$body = $newSlide.Shapes.Item(2)

# Position/size override matching the authored slide.
$body.Left = 1371600 / 12700.0
$body.Top = 2286000 / 12700.0
$body.Width = 9601200 / 12700.0
$body.Height = 4364182 / 12700.0

$bodyTr = $body.TextFrame.TextRange
$bodyTr.Text = "The project is developed using Python and the python can be downloaded from the below GitHub Link:`r`r" + `
    "https://github.com/kgopal1982/Analytics/blob/master/HomeLoanAnalytics/HomeLoanAssessment.ipynb`r`r`r`r"

# Paragraph 3 holds the GitHub URL; turn it into a hyperlink.
$linkPara = $bodyTr.Paragraphs(3, 1)
$linkAs = $linkPara.ActionSettings.Item(1)
$linkAs.Hyperlink.Address = "https://github.com/kgopal1982/Analytics/blob/master/HomeLoanAnalytics/HomeLoanAssessment.ipynb"

# Paragraphs 3, 4, 5 and 7 are flush-left with no bullet in the authored
# slide; paragraphs 2 and 6 keep the placeholder's default (bulleted)
# formatting since they were left untouched (and stay blank).
$bodyTr.Paragraphs(3, 1).ParagraphFormat.Bullet.Visible = 0
$bodyTr.Paragraphs(4, 1).ParagraphFormat.Bullet.Visible = 0
$bodyTr.Paragraphs(5, 1).ParagraphFormat.Bullet.Visible = 0
$bodyTr.Paragraphs(7, 1).ParagraphFormat.Bullet.Visible = 0

# Paragraphs 2 and 6 are plain blank lines; re-clear their text so they
# collapse back down to bare empty paragraphs instead of keeping a
# leftover empty run.
$bodyTr.Paragraphs(2, 1).Text = ""
$bodyTr.Paragraphs(6, 1).Text = ""
